$wb = $excel.ActiveWorkbook

    $wb.Worksheets.Item(1).Name = "summ31627255"
    $wb.Worksheets.Item(2).Name = "summ31721189"
    $wb.Worksheets.Item(3).Name = "summ31807098"
    $wb.Worksheets.Item(4).Name = "summ31891869"
    $wb.Worksheets.Item(5).Name = "summ31986236"
    $wb.Worksheets.Item(6).Name = "summ32068371"
    $wb.Worksheets.Item(7).Name = "summ32154618"
    $wb.Worksheets.Item(8).Name = "summ32234003"
    $wb.Worksheets.Item(9).Name = "summ32315394"
    $wb.Worksheets.Item(10).Name = "summ32390891"
    $wb.Worksheets.Item(11).Name = "summ32486227"
    $wb.Worksheets.Item(12).Name = "summ32569554"
    $wb.Worksheets.Item(13).Name = "summ32647157"
    $wb.Worksheets.Item(14).Name = "summ32730740"
    $wb.Worksheets.Item(15).Name = "summ32805628"
    $wb.Worksheets.Item(16).Name = "summ32896329"
    $wb.Worksheets.Item(17).Name = "summ32982273"
    $wb.Worksheets.Item(18).Name = "summ33056340"
    $wb.Worksheets.Item(19).Name = "summ33148114"
    $wb.Worksheets.Item(20).Name = "summ33237158"
    $wb.Worksheets.Item(21).Name = "summ33350079"
    $wb.Worksheets.Item(22).Name = "summ33488701"
    $wb.Worksheets.Item(23).Name = "summ33588053"
    $wb.Worksheets.Item(24).Name = "summ33720387"
    $wb.Worksheets.Item(25).Name = "summ33916581"
    $wb.Worksheets.Item(26).Name = "summ34189833"
    $wb.Worksheets.Item(27).Name = "summ34346234"
    $wb.Worksheets.Item(28).Name = "summ34474443"
    $wb.Worksheets.Item(29).Name = "summ34626098"
    $wb.Worksheets.Item(30).Name = "summ34774605"
    $wb.Worksheets.Item(31).Name = "summ34906487"
    $wb.Worksheets.Item(32).Name = "summ34987363"
    $wb.Worksheets.Item(33).Name = "summ35087362"
    $wb.Worksheets.Item(34).Name = "summ35178360"
    $wb.Worksheets.Item(35).Name = "summ35270637"
    $wb.Worksheets.Item(36).Name = "summ35355014"
    $wb.Worksheets.Item(37).Name = "summ35437752"
    $wb.Worksheets.Item(38).Name = "summ35537736"
    $wb.Worksheets.Item(39).Name = "summ35621164"
    $wb.Worksheets.Item(40).Name = "summ35719566"
    $wb.Worksheets.Item(41).Name = "summ35810983"
    $wb.Worksheets.Item(42).Name = "summ35903917"
    $wb.Worksheets.Item(43).Name = "summ35988156"
    $wb.Worksheets.Item(44).Name = "summ36072613"
    $wb.Worksheets.Item(45).Name = "summ36156814"
    $wb.Worksheets.Item(46).Name = "summ36253902"
    $wb.Worksheets.Item(47).Name = "summ36337004"
    $wb.Worksheets.Item(48).Name = "summ36425728"
    $wb.Worksheets.Item(49).Name = "summ36506609"
    $wb.Worksheets.Item(50).Name = "summ36589304"
